# Update Name of Algo
# Apply corrected values to the KNN imputation result data (columns A and B)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 5.697000000000001
$ws.Range("A8").Value = -22.366
$ws.Range("A10").Value = -21.604
$ws.Range("A12").Value = -21.606
$ws.Range("B15").Value = 4.867
$ws.Range("A18").Value = -21.313
$ws.Range("B18").Value = 7.007
$ws.Range("B20").Value = 6.406999999999999
$ws.Range("B29").Value = 5.645
$ws.Range("B30").Value = 6.208
$ws.Range("B31").Value = 6.367000000000001
$ws.Range("A37").Value = -19.92
$ws.Range("B40").Value = 8.962
$ws.Range("B50").Value = 4.807
$ws.Range("A55").Value = -21.797
$ws.Range("A68").Value = -21.507
$ws.Range("B68").Value = 5.881
$ws.Range("B76").Value = 6.343000000000001
$ws.Range("A77").Value = -20.637
$ws.Range("A78").Value = -20.082
$ws.Range("A81").Value = -21.818
$ws.Range("A82").Value = -22.077
$ws.Range("B87").Value = 4.636
$ws.Range("B88").Value = 4.858000000000001
$ws.Range("B96").Value = 6.692
$ws.Range("B98").Value = 5.355
$ws.Range("B101").Value = 7.782000000000001
$ws.Range("B102").Value = 7.747
